$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column at N (bool_value), shifting old N:P -> O:Q
$ws.Columns("N").Insert()

# Give the new column the same width as column M so they render as one
# visually-consistent block (matches width of the "title" column).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# 2. New header for the inserted column
$ws.Range("N1").Value = "bool_value"

# 3. Fix mislabeled attribute name: row 126 should say "genotype" (singular),
#    matching the rest of that row (E126/G126 already say "genotype").
$ws.Range("B126").Value = "genotype"

# 4. Fix the swapped "maximum values" / "order" flags between rows 59 and 60:
#    row 59 incorrectly carried an order index of 5 and a false "display" flag,
#    while row 60 incorrectly carried the true flag without the order index.
$ws.Range("I59").ClearContents()
$ws.Range("J59").Value = $true

$ws.Range("I60").Value = 5
$ws.Range("J60").Value = $false

# 5. Update the hidden _FilterDatabase defined name so it covers the new column.
$fd = $wb.Names.Item("airr_schema_defs!_FilterDatabase")
$fd.RefersTo = "=airr_schema_defs!`$A`$1:`$Q`$144"

# 6. Restore the active selection to D5 (was D15).
$ws.Range("D5").Select() | Out-Null
